# Update the (encrypted) values shown in column B for rows 3, 14 and 15,
# and move the sheet's active-cell selection to B8.
#
# Order matters: new literal strings are appended to the shared-string
# table in the order they're first written, and the target file expects
# row 15's new string to land at index 29 and row 3's new string at index
# 30 - so B15 must be written before B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "U2FsdGVkX1/7B7N32aaq9jrYalUe5yoZD7xuJE34Zdpb845TlihaGy7rZTWYz3JI6vibM4VaDWoouTDS00TEL4QatZIMXYGLQCWnqakX1xjzsNu/txXJ/+ETeJs1ZP9B4dzIrmkZg9w+jZUXgj3YqPL80LWMmWKBrqXc4qUUAR6mHtI6oawguolHaKQOoOnuxvoz5jyMDBj4Hp0GuWm7c3Bhx5tUQ8BvPHICVGgBtt4mVLTiXwiVhiEAVqMBkb9OfS/j8HfeXQwgGeqlsYCFDw=="

$ws.Range("B3").Value = "U2FsdGVkX19PQST6oRb3uUB6cyFGdsIe3QUaut36xRhgLrSBER105qGmrDnJXxbaWymmB+ThxbaUB/2O+Ds43dTIknLYob/k0M0RbX2MDgrctlkZFdFNYVqlWnLM3by+hv3UN1QizejtpYiLVxV7/ArnVBXUjSX2aIKHQnUvA+spseey7QSREKnbVdtxq41v+qUg9R26IAsgSlcxprgipHlJFeqWr5YUz04TPmqAgfkEDM3TDF/c85UlZSjbRuIdhmCL3UsET/2dGAIUhw8t0zUEPVlDGc1aGLYOwX2OEk00ADNWWXrQFObHiQg6yQmcIE6XkjUB5naO7fSZC/jolVkqQvGMS5J/LLqTxJfTIh4="

$ws.Range("B14").Value = "U2FsdGVkX1+lUCchG5dRfcxw4S6OG26+o0pLVrsHsFW6izTIrGzFN3EKcCOObQdHnIi0fHt5TvjS25V5vkunwXgup/KU8pe2ncLBHfpEe6g4/nYIZz+davjGU9E1GqCyfLmGpPbSIOb6Us/NUpWFPibhNbsn7w4ROQ2ylZ3AcWItE5Sy9zlabp1ix2lNiUfz/qLAIGZKNi+YdkCqJ9oDYO093kPdwFKbgzpiiui1gzg9ZcjKZ5WJ/BPZ/vd5CPtJmayjAWVKsW1U7cXelWuOuQjRajUvAn2y3KEUfuqWKVLnFMAaG3+w8Hx9+XaYfmE2TjRPjX2RvMycforapKZSwtFEelmItsdZIniZG3rZO0H6Ozm64A/8EBy8PBLqyw7gYLNOislnQTJiaWEmQViAN+M+RJGx8zp5EXfrRkvRnpcpDQc59j4rB7qaEZyVEDcjDwWocqqyPnXYPjBSfbyRIAxl3jocGxpZ5nwfLeQgt2GTIrNRl2YR3/kz7/ak7mG28QoFfTkPaGrbKBrM3tCBvRcFXZftifV45linnLbqSgwymgVrvO8Pf92AOSHM4Aewdvag0V3mrUCcVkjv7PunW8fOU1YCv4hX2XR5WKYk6J0FrW5/0TsDtj3vqXFIpxSUdlwkk9brn4XxZCPd/hGdYH9GjF9HCNLf1fgfH6X9SPFZ2pnmo7ni76IhpjlG1cKw2bYI7TtcpikDEbFCNpuQxy27cxfGJJ2HzA2wHzbz17q2XaGuCoxi9/4oYcepARlGoYpZl/Bf+DBwXMTbeFZ16vkoZE433htfx9xcZfIf5sDR5xmwzq70olbauMilQsRW3MaiqI8jPnn8SVv/kyh5seIKBVlLgvEkahn+d02GejGiCJN4HOjYXqehA8D8t8sUnLJblsjqgQc4uzBEukntgB1AZvQoUHsJHJAGm7Lp9LgTwvys4SPneJLWIxZ2Dnh/"

$ws.Range("B8").Select()
